$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14/15 swap: Polkadot/Polygon -> Polygon/Polkadot (with new values) ---
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.797"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.04%  "

# --- Row 42/43 swap: VeChain/InjectiveProtocol -> InjectiveProtocol/VeChain (with new values) ---
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.02%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.58%  "

# --- Remaining price/volume updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.030.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("E6").Value = "  +1.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.23"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.11%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "

$ws.Range("E10").Value = "  +3.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.349.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.054.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.045.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +17.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0894"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0615"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0886"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.37%  "

$ws.Range("E38").Value = "  +3.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +26.43%  "

$ws.Range("E44").Value = "  -1.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("E46").Value = "  +5.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.284.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.239.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.51%  "

